$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-5 from 45175 to 45183
$ws.Range("C2:C5").Value = 45183
